# Fixed naive component forecaster bug - Presentation state 11.02.
#
# The "naive QoQ error" matrix gets a new Q0 (column B) error series.
# For row 2 (the very first vintage) nothing is inserted - it only loses
# its trailing K2 value. For every other populated row (3..24) a brand
# new near-zero "Q0" residual is inserted at column B, every existing
# value shifts one column to the right, and whatever value would have
# spilled past column K (the last column of the triangle) is dropped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column-B ("Q0") values keyed by row number.
$newB = @{
    3  = [double]"1.052939957446597E-10"
    4  = [double]"2.480007310623478E-10"
    5  = [double]"1.503854694107076E-07"
    6  = [double]"-3.965936795080616E-07"
    7  = [double]"-3.930720193778825E-10"
    8  = [double]"1.181302580199883E-07"
    9  = [double]"2.965444589886346E-07"
    10 = [double]"3.593882045849206E-07"
    11 = [double]"-1.035781544145298E-07"
    12 = [double]"-3.87512216759589E-10"
    13 = [double]"-1.07388789361007E-07"
    14 = [double]"-1.035472805832605E-07"
    15 = [double]"6.303355340908645E-06"
    16 = [double]"-2.375649628613696E-07"
    17 = [double]"3.720025918141356E-07"
    18 = [double]"3.829984367986761E-07"
    19 = [double]"-3.160475492397508E-06"
    20 = [double]"-4.101096154340844E-08"
    21 = [double]"-1.831659499074156E-07"
    22 = [double]"2.770877186031306E-07"
    23 = [double]"2.29775004800814E-07"
    24 = [double]"-1.554241066958895E-07"
}

# Last populated data column (1-indexed; A=1, B=2, ... K=11) for each row,
# taken from the sheet BEFORE this edit is applied.
$lastCol = @{
    2  = 11
    3  = 11
    4  = 11
    5  = 11
    6  = 11
    7  = 11
    8  = 11
    9  = 11
    10 = 11
    11 = 11
    12 = 11
    13 = 11
    14 = 11
    15 = 10
    16 = 9
    17 = 8
    18 = 7
    19 = 6
    20 = 5
    21 = 4
    22 = 3
    23 = 2
    24 = 1
}

# Row 2 is untouched apart from dropping its last (K2) value - no new
# Q0 column was inserted for this vintage.
$ws.Cells.Item(2, 11).ClearContents()

# Rows 3..24: shift existing values one column to the right (starting
# from the rightmost populated column and working back to C so we never
# clobber a value before it has been copied), dropping anything that
# would spill past column K, then write the new Q0 value into column B.
for ($row = 3; $row -le 24; $row++) {
    $last = $lastCol[$row]
    $target = [Math]::Min($last + 1, 11)

    for ($col = $target; $col -ge 3; $col--) {
        $srcVal = $ws.Cells.Item($row, $col - 1).Value2
        if ($srcVal -eq $null) {
            $ws.Cells.Item($row, $col).ClearContents()
        } else {
            $ws.Cells.Item($row, $col).Value = $srcVal
        }
    }

    $ws.Cells.Item($row, 2).Value = $newB[$row]
}

Write-Output "done"
